$d = $word.ActiveDocument

# Insert a brand-new paragraph before the document's current first
# paragraph (the "En mi experiencia laboral..." text) to hold the new
# centered, bold title line.
$firstPara = $d.Paragraphs.Item(1)
$insertPoint = $firstPara.Range
$insertPoint.Collapse(1)  # wdCollapseStart
$insertPoint.InsertParagraphBefore()

# The freshly inserted paragraph is now Paragraphs(1); set its text.
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Text = "IMPORTANCIA DE LAS HABILIDADES BLANDAS EN EL ENTORNO LABORAL."

# Re-acquire the paragraph/range after the text write, then apply the
# paragraph + character formatting: centered alignment, bold, black
# color, 13.5pt (sz/szCs = 27 half-points) text.
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range

$titlePara.Format.Alignment = 1  # wdAlignParagraphCenter

$titleRange.Font.Bold = $true
$titleRange.Font.BoldBi = $true
$titleRange.Font.Color = 0        # wdColorAutomatic/black (RGB 0,0,0)
$titleRange.Font.Size = 13.5
$titleRange.Font.SizeBi = 13.5
